# Update imputed values in the RandomForest result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Cell = "A3";  Value = -22.08960000000001 },
    @{ Cell = "A21"; Value = -19.90869999999998 },
    @{ Cell = "A23"; Value = -20.58059999999999 },
    @{ Cell = "A25"; Value = -21.51449999999999 },
    @{ Cell = "D27"; Value = -8.804300000000003 },
    @{ Cell = "D31"; Value = -8.231000000000003 },
    @{ Cell = "D39"; Value = -8.116899999999998 },
    @{ Cell = "D48"; Value = -7.352599999999998 },
    @{ Cell = "D51"; Value = -7.853700000000003 },
    @{ Cell = "D52"; Value = -7.6518 },
    @{ Cell = "A53"; Value = -21.8154 },
    @{ Cell = "D55"; Value = -8.251899999999997 },
    @{ Cell = "D56"; Value = -7.940999999999997 },
    @{ Cell = "A57"; Value = -22.1693 },
    @{ Cell = "D57"; Value = -8.567100000000005 },
    @{ Cell = "A59"; Value = -22.19179999999999 },
    @{ Cell = "A69"; Value = -21.6352 },
    @{ Cell = "D73"; Value = -7.825099999999996 },
    @{ Cell = "A79"; Value = -20.46500000000001 },
    @{ Cell = "A83"; Value = -21.9217 },
    @{ Cell = "D89"; Value = -5.786900000000002 },
    @{ Cell = "D90"; Value = -8.035000000000004 },
    @{ Cell = "A93"; Value = -21.3802 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
